$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.06"
$ws.Range("E2").Value = "'-1.08%"
$ws.Range("D3").Value = "'26.65"
$ws.Range("E3").Value = "'-2.90%"
$ws.Range("D4").Value = "'4.877"
$ws.Range("E4").Value = "'0.93%"
$ws.Range("D5").Value = "'0.06328"
$ws.Range("E5").Value = "'-0.36%"
$ws.Range("D6").Value = "'6.915"
$ws.Range("E6").Value = "'-0.38%"
$ws.Range("D7").Value = "'3.319"
$ws.Range("E7").Value = "'-1.97%"
$ws.Range("E8").Value = "'35.01%"
$ws.Range("D9").Value = "'0.8748"
$ws.Range("E9").Value = "'-0.49%"
$ws.Range("D10").Value = "'0.1574"
$ws.Range("E10").Value = "'7.06%"
$ws.Range("D11").Value = "'0.05011"
$ws.Range("E11").Value = "'-2.82%"
$ws.Range("D12").Value = "'0.07484"
$ws.Range("E12").Value = "'2.22%"
$ws.Range("D13").Value = "'0.02955"
$ws.Range("E13").Value = "'-5.84%"
$ws.Range("D14").Value = "'0.09061"
$ws.Range("E14").Value = "'-0.07%"
$ws.Range("D15").Value = "'0.001584"
$ws.Range("E15").Value = "'1.23%"
$ws.Range("D16").Value = "'0.0006345"
$ws.Range("E16").Value = "'0.86%"
$ws.Range("D17").Value = "'0.006013"
$ws.Range("E17").Value = "'-0.13%"
$ws.Range("D18").Value = "'3.446"
$ws.Range("E18").Value = "'0.05%"
$ws.Range("E19").Value = "'-0.36%"
$ws.Range("D21").Value = "'0.1334"
$ws.Range("E21").Value = "'1.75%"
$ws.Range("D22").Value = "'3.905"
$ws.Range("E22").Value = "'1.39%"
$ws.Range("D23").Value = "'0.04368"
$ws.Range("E23").Value = "'0.96%"
$ws.Range("E24").Value = "'-0.55%"
$ws.Range("E25").Value = "'-2.04%"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("D27").Value = "'0.0001616"
$ws.Range("E27").Value = "'-4.37%"
$ws.Range("D40").Value = "'0.04100"
$ws.Range("E40").Value = "'0.23%"
$ws.Range("D41").Value = "'0.007020"
$ws.Range("E41").Value = "'5.06%"
$ws.Range("D42").Value = "'0.1171"
$ws.Range("E42").Value = "'0.53%"
$ws.Range("D43").Value = "'0.002129"
$ws.Range("E43").Value = "'-3.18%"
$ws.Range("D44").Value = "'0.01080"
$ws.Range("E44").Value = "'-18.25%"
$ws.Range("D45").Value = "'0.00005304"
$ws.Range("E45").Value = "'1.71%"
$ws.Range("D46").Value = "'0.02000"
$ws.Range("E46").Value = "'-11.06%"
$ws.Range("D47").Value = "'1.486"
$ws.Range("E47").Value = "'-37.50%"
